# Shift the "hour" (column B) and "Offerte" (column C) series up by one
# row so the real-time line chart starts counting from hour 8, and add
# the newly-arrived data points (C12 = 201, B18 = 24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bValues = @(8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24)
$cValues = @(0, 52, 103, 140, 130, 128, 111, 90, 106, $null, 201, $null, $null, $null, $null, $null, $null)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = 2 + $i

    if ($null -eq $bValues[$i]) {
        $ws.Cells.Item($row, 2).Value = ""
    } else {
        $ws.Cells.Item($row, 2).Value = $bValues[$i]
    }

    if ($null -eq $cValues[$i]) {
        $ws.Cells.Item($row, 3).Value = ""
    } else {
        $ws.Cells.Item($row, 3).Value = $cValues[$i]
    }
}
